$d = $word.ActiveDocument
$s = $d.Styles.Item(11)
Write-Output ("BuiltIn=" + $s.BuiltIn)
Write-Output ("Priority=" + $s.Priority)
Write-Output ("QuickStyle=" + $s.QuickStyle)
Write-Output ("Locked=" + $s.Locked)
